# Slide 14 ("INICIALIZACAO DA FERRAMENTA") has two click-triggered entrance
# (fade) animations that target the picture "Imagem 6" (shape id=7) and the
# shape "Seta para a direita 1" (shape id=2). The edit swaps which of the
# two click steps animates which shape (the picture now appears on the
# earlier click, the arrow on the later one) while leaving every other
# animation on the slide untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(14)
$seq = $s.TimeLine.MainSequence

# Shapes involved in the swap (indices are positions in Slide.Shapes).
$shapeImagem6 = $s.Shapes.Item(5)   # id=7 "Imagem 6"
$shapeSeta1   = $s.Shapes.Item(6)   # id=2 "Seta para a direita 1"
$shapeSeta7   = $s.Shapes.Item(7)   # id=8 "Seta para a direita 7"

# The main click sequence on this slide is currently:
#   1: Espaco Reservado para Texto 3 (spid=4)
#   2: Imagem 5                       (spid=6)
#   3: Imagem 6                       (spid=7)
#   4: Seta para a direita 1          (spid=2)
#   5: Seta para a direita 7          (spid=8)
# and it needs to become:
#   1: Espaco Reservado para Texto 3 (spid=4)   -- unchanged
#   2: Imagem 5                       (spid=6)   -- unchanged
#   3: Seta para a direita 1          (spid=2)
#   4: Imagem 6                       (spid=7)
#   5: Seta para a direita 7          (spid=8)   -- unchanged (re-appended as-is)
#
# New effects are always appended at the end of the sequence, so the three
# trailing effects are removed and re-added in the desired order.
$seq.Item(4).Delete()   # remove "Seta para a direita 1" effect (was position 4)
$seq.Item(3).Delete()   # remove "Imagem 6" effect (was position 3)
$seq.Item(3).Delete()   # remove "Seta para a direita 7" effect (now at position 3)

$null = $seq.AddEffect($shapeSeta1, 10)    # -> position 3
$null = $seq.AddEffect($shapeImagem6, 10)  # -> position 4
$null = $seq.AddEffect($shapeSeta7, 10)    # -> position 5 (restored)
